$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 1062.458
$ws.Range("E2").Value = 430065000.0
$ws.Range("F2").Value = 946938.6000000001
$ws.Range("G2").Value = 516873.6
$ws.Range("H2").Value = 125.61599999999999
$ws.Range("I2").Value = 27
$ws.Range("J2").Value = 146.17199999999997
$ws.Range("K2").Value = 20.555
$ws.Range("L2").Value = 574.304
$ws.Range("M2").Value = 488710.22727272724
$ws.Range("N2").Value = 531.9879775280899
